$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A23").Value = "floor8_ptsw1"
$ws.Range("B24").Value = "192.168.0.3"
$ws.Range("A25").Value = "asher1"
$ws.Range("B25").Value = "the men1"
